$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A7").Value = "School 1"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Surname 1"
$ws.Range("D7").Value = "Surname 1"
$ws.Range("E7").Value = "6018131X"
$ws.Range("F7").Value = 43647
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("G7").Value = "U"

[void]$ws.Range("B12").Select()
